$d = $word.ActiveDocument

# "away time" -> "last seen time" in the bulleted feature list under
# "user info".
$d.Content.Find.Execute("away time", $false, $false, $false, $false, `
    $false, $true, 1, $false, "last seen time", 2)

# Word's "_GoBack" bookmark is a singleton that always marks the location
# of the most recent edit -- adding it here automatically removes it from
# wherever it previously sat (after "chat synchronization, and chat
# messages" in the multicast-channel paragraph).
$rng = $d.Content
$rng.Find.Execute("last seen")
$editPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $editPoint)
